$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 7356.174
$ws.Cells.Item(116, 9).Value = 7419.6
$ws.Cells.Item(116, 11).Value = 7419.6
$ws.Cells.Item(116, 13).Value = -3977.6
$ws.Cells.Item(129, 8).Value = 2800.1
$ws.Cells.Item(129, 9).Value = 2348.25
$ws.Cells.Item(129, 10).Value = 3101.3333
$ws.Cells.Item(129, 11).Value = 7044.75
$ws.Cells.Item(129, 12).Value = 9303.999899999999
$ws.Cells.Item(129, 13).Value = -2044.75
$ws.Cells.Item(129, 14).Value = -19303.9999
$ws.Cells.Item(132, 8).Value = 2235.8462
$ws.Cells.Item(132, 9).Value = 2235.8462
$ws.Cells.Item(132, 11).Value = 6707.5386
$ws.Cells.Item(132, 13).Value = -4177.5386
$ws.Cells.Item(135, 8).Value = 694.2381
$ws.Cells.Item(135, 9).Value = 628.9
$ws.Cells.Item(135, 10).Value = 2001
$ws.Cells.Item(135, 11).Value = 5660.099999999999
$ws.Cells.Item(135, 12).Value = 18009
$ws.Cells.Item(135, 13).Value = -3125.099999999999
$ws.Cells.Item(135, 14).Value = -23079
$ws.Cells.Item(138, 8).Value = 2848.7437
$ws.Cells.Item(138, 9).Value = 2393.4707
$ws.Cells.Item(138, 10).Value = 3200.5454
$ws.Cells.Item(138, 11).Value = 7180.4121
$ws.Cells.Item(138, 12).Value = 9601.636200000001
$ws.Cells.Item(138, 13).Value = -2040.4121
$ws.Cells.Item(138, 14).Value = -19881.6362

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4754.4873
$ws.Cells.Item(32, 9).Value = 3961.7058
$ws.Cells.Item(32, 11).Value = 3961.7058
$ws.Cells.Item(32, 13).Value = -3674.7058
$ws.Cells.Item(45, 8).Value = 62504450
$ws.Cells.Item(45, 9).Value = 111112530
$ws.Cells.Item(45, 11).Value = 111112530
$ws.Cells.Item(45, 13).Value = -111112153
$ws.Cells.Item(74, 8).Value = 11495826
$ws.Cells.Item(74, 9).Value = 13334901
$ws.Cells.Item(74, 11).Value = 13334901
$ws.Cells.Item(74, 13).Value = -13334027
$ws.Cells.Item(77, 8).Value = 11495826
$ws.Cells.Item(77, 9).Value = 13334901
$ws.Cells.Item(77, 11).Value = 66674505
$ws.Cells.Item(77, 13).Value = -66670137
$ws.Cells.Item(122, 8).Value = 2544.5833
$ws.Cells.Item(122, 9).Value = 1866.8182
$ws.Cells.Item(122, 11).Value = 5600.4546
$ws.Cells.Item(122, 13).Value = -3150.4546

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2888.7334
$ws.Cells.Item(20, 10).Value = 3686
$ws.Cells.Item(20, 12).Value = 3686
$ws.Cells.Item(20, 14).Value = -4180
$ws.Cells.Item(62, 8).Value = 49000
$ws.Cells.Item(62, 10).Value = 49000
$ws.Cells.Item(62, 12).Value = 49000
$ws.Cells.Item(62, 14).Value = -50372
$ws.Cells.Item(65, 8).Value = 49000
$ws.Cells.Item(65, 10).Value = 49000
$ws.Cells.Item(65, 12).Value = 147000
$ws.Cells.Item(65, 14).Value = -153864
$ws.Cells.Item(74, 8).Value = 59998.5
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).Value = ""
$ws.Cells.Item(75, 8).Value = 40000
$ws.Cells.Item(75, 9).Value = 10000
$ws.Cells.Item(75, 10).Value = 70000
$ws.Cells.Item(75, 11).Value = 10000
$ws.Cells.Item(75, 12).Value = 70000
$ws.Cells.Item(75, 13).Value = -9064
$ws.Cells.Item(75, 14).Value = -71872
$ws.Cells.Item(77, 8).Value = 59998.5
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).Value = ""
$ws.Cells.Item(78, 8).Value = 40000
$ws.Cells.Item(78, 9).Value = 10000
$ws.Cells.Item(78, 10).Value = 70000
$ws.Cells.Item(78, 11).Value = 30000
$ws.Cells.Item(78, 12).Value = 210000
$ws.Cells.Item(78, 13).Value = -25320
$ws.Cells.Item(78, 14).Value = -219360
$ws.Cells.Item(96, 8).Value = 28459.5
$ws.Cells.Item(96, 9).Value = 28459.5
$ws.Cells.Item(96, 11).Value = 28459.5
$ws.Cells.Item(96, 13).Value = -25713.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(36, 8).Value = 19450
$ws.Cells.Item(36, 10).Value = 19450
$ws.Cells.Item(36, 12).Value = 19450
$ws.Cells.Item(36, 14).Value = -20226
$ws.Cells.Item(38, 8).Value = 15000
$ws.Cells.Item(38, 10).Value = 15000
$ws.Cells.Item(38, 12).Value = 15000
$ws.Cells.Item(38, 14).Value = -15754
$ws.Cells.Item(39, 8).Value = 5000
$ws.Cells.Item(39, 9).Value = 5000
$ws.Cells.Item(39, 11).Value = 5000
$ws.Cells.Item(39, 13).Value = -4609
$ws.Cells.Item(40, 8).Value = 19450
$ws.Cells.Item(40, 10).Value = 19450
$ws.Cells.Item(40, 12).Value = 19450
$ws.Cells.Item(40, 14).Value = -19770
$ws.Cells.Item(46, 8).Value = 15000
$ws.Cells.Item(46, 10).Value = 15000
$ws.Cells.Item(46, 12).Value = 15000
$ws.Cells.Item(46, 14).Value = -15422
$ws.Cells.Item(49, 8).Value = 5000
$ws.Cells.Item(49, 9).Value = 5000
$ws.Cells.Item(49, 11).Value = 5000
$ws.Cells.Item(49, 13).Value = -4818

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 6476205
$ws.Cells.Item(4, 9).Value = 9850256
$ws.Cells.Item(4, 11).Value = 29550768
$ws.Cells.Item(4, 13).Value = -29550656
$ws.Cells.Item(25, 8).Value = 1308.4667
$ws.Cells.Item(25, 10).Value = 1495.2307
$ws.Cells.Item(25, 12).Value = 4485.6921
$ws.Cells.Item(25, 14).Value = -4823.6921
$ws.Cells.Item(30, 8).Value = 1308.4667
$ws.Cells.Item(30, 10).Value = 1495.2307
$ws.Cells.Item(30, 12).Value = 4485.6921
$ws.Cells.Item(30, 14).Value = -4689.6921
$ws.Cells.Item(86, 8).Value = 2519.5715
$ws.Cells.Item(86, 10).Value = 4149.25
$ws.Cells.Item(86, 12).Value = 12447.75
$ws.Cells.Item(86, 14).Value = -14819.75
$ws.Cells.Item(89, 8).Value = 2519.5715
$ws.Cells.Item(89, 10).Value = 4149.25
$ws.Cells.Item(89, 12).Value = 37343.25
$ws.Cells.Item(89, 14).Value = -49199.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 16772
$ws.Cells.Item(70, 9).Value = 5359.4614
$ws.Cells.Item(70, 10).Value = 41499.168
$ws.Cells.Item(70, 11).Value = 5359.4614
$ws.Cells.Item(70, 12).Value = 41499.168
$ws.Cells.Item(70, 13).Value = -5089.4614
$ws.Cells.Item(70, 14).Value = -42039.168
$ws.Cells.Item(73, 8).Value = 16772
$ws.Cells.Item(73, 9).Value = 5359.4614
$ws.Cells.Item(73, 10).Value = 41499.168
$ws.Cells.Item(73, 11).Value = 5359.4614
$ws.Cells.Item(73, 12).Value = 41499.168
$ws.Cells.Item(73, 13).Value = -4423.4614
$ws.Cells.Item(73, 14).Value = -43371.168
$ws.Cells.Item(102, 8).Value = 2809.1765
$ws.Cells.Item(102, 9).Value = 1869.5333
$ws.Cells.Item(102, 10).Value = 9856.5
$ws.Cells.Item(102, 11).Value = 1869.5333
$ws.Cells.Item(102, 12).Value = 9856.5
$ws.Cells.Item(102, 13).Value = -247.5333000000001
$ws.Cells.Item(102, 14).Value = -13100.5
$ws.Cells.Item(126, 8).Value = 2974.5264
$ws.Cells.Item(126, 9).Value = 2000.1333
$ws.Cells.Item(126, 10).Value = 6628.5
$ws.Cells.Item(126, 11).Value = 6000.3999
$ws.Cells.Item(126, 12).Value = 19885.5
$ws.Cells.Item(126, 13).Value = -3530.3999
$ws.Cells.Item(126, 14).Value = -24825.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6183.9165
$ws.Cells.Item(40, 9).Value = 5385.2856
$ws.Cells.Item(40, 11).Value = 5385.2856
$ws.Cells.Item(40, 13).Value = -5249.2856
$ws.Cells.Item(122, 8).Value = 10494.875
$ws.Cells.Item(122, 9).Value = 8790.799999999999
$ws.Cells.Item(122, 11).Value = 26372.4
$ws.Cells.Item(122, 13).Value = -23922.4
$ws.Cells.Item(128, 8).Value = 112465
$ws.Cells.Item(128, 10).Value = 112465
$ws.Cells.Item(128, 12).Value = 112465
$ws.Cells.Item(128, 14).Value = -122425

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48, 8).Value = 20000
$ws.Cells.Item(48, 9).Value = 20000
$ws.Cells.Item(48, 11).Value = 20000
$ws.Cells.Item(48, 13).Value = -19431
$ws.Cells.Item(50, 8).Value = 22500
$ws.Cells.Item(50, 10).Value = 22500
$ws.Cells.Item(50, 12).Value = 22500
$ws.Cells.Item(50, 14).Value = -23762
$ws.Cells.Item(107, 8).Value = 382.34616
$ws.Cells.Item(107, 9).Value = 344.61905
$ws.Cells.Item(107, 11).Value = 1033.85715
$ws.Cells.Item(107, 13).Value = 886.14285
$ws.Cells.Item(122, 8).Value = 8593.038
$ws.Cells.Item(122, 9).Value = 2150.875
$ws.Cells.Item(122, 11).Value = 6452.625
$ws.Cells.Item(122, 13).Value = -4002.625
